$d = $word.ActiveDocument

# Locate the "Tortoise Media" text inside the sentence about the seven indicators.
$rng = $d.Content
$found = $rng.Find.Execute("Tortoise Media", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)

# $rng now collapses to the found text ("Tortoise Media"); turn it into a hyperlink.
$h = $d.Hyperlinks.Add($rng, "https://www.tortoisemedia.com/intelligence/global-ai/", "", "", "Tortoise Media")
$h.Range.Font.Name = "Times New Roman"
$h.Range.Font.NameAscii = "Times New Roman"
$h.Range.Font.NameBi = "Times New Roman"
